$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new Portuguese-language scenario value in the row below the
# existing ones (A3 = "joiceAbreu", A4 = "JoiceGGG", A5 = "JOOCEAA").
$rng = $ws.Range("A5")
$rng.Value = "JOOCEAA"

# Give it the same kind of "top + bottom" thin border used for the other
# login/senha rows above it.
$rng.Borders.Item(8).LineStyle = 1
$rng.Borders.Item(9).LineStyle = 1

# Move/save the active selection onto the newly entered cell.
$rng.Select()
